$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Could not find text: $old"
    }
}

Replace-Text "WHAT IS THE SMARTCARD?" "Что такое SmartCard?"

Replace-Text "The SmartCard is a way to hold SMART on a physical card. You don’t need a phone to hold your SMART. If you have ever had a problem with using a phone inside a building…you will get why this is ideal. Only the merchant needs the internet connection. As simple to use as a credit or debit card but not controlled by any 3rd party processor. And not a preloaded card using a visa fiat system…this is real crypto to crypto." "SmartCard — это способ хранить и совершать платежи в SMART, используя физическую карту. Вам не нужен даже телефон. Вероятно, вы сталкивались с проблемой качества мобильного интернета, поэтому вы поймёте наше стремление пойти дальше. Теперь только продавцы должны иметь интернет-соединение. Это так же просто, как использование дебетовой или кредитной карты, но происходит это без посредничества VISA или любых других платежных процессоров. Поэтому транзакции не контролируются третьими сторонами. Это настоящие платежи в крипто."

Replace-Text "FUTURE POTENTIAL OF SMARTCARD…" "Будущий потенциал SmartCard"

Replace-Text "Debit and credit card payments have arguably become the dominant way for many to carry out most, if not all, of their daily transactions. Adoption of card payments became increasingly simple for smaller businesses, and now quick and affordable smartphone-powered solutions are already a reality." "Дебетовые и кредитные карты заняли доминирующие позиции в сегменте повседневных транзакций для многих людей. Принятие платежей посредством карт, оплата через смартфоны или иные умные устройства — быстрое и удобное решение. Оно идеально в случае индивидуальных предпринимателей или малого бизнеса."

Replace-Text "SmartCard$([char]0x00A0)emerged offering to solve many of the problems facing traditional systems.$([char]0x00A0)SmartCard$([char]0x00A0) is blockchain-based solution that aims to improve on the digital payment card model. Our main goals are to streamline crypto transactions and make them practical for use in real-life payment scenarios. At the same time, we aim to drastically reduce fees and confirmation time frames over those found in current solutions." "SmartCard ставит перед собой задачу избавиться от многих недостатков, которые есть у традиционных систем. SmartCard — это основанное на базе Blockchain решение, призванное улучшить модель оплаты с помощью карт. Наша главная цель — навсегда изменить криптовалютные транзакции, сделать их использование в повседневной жизни простым и удобным. В тоже время, мы стремимся иметь наименьшую комиссию и время подтверждения — для вашего удобства."

Replace-Text "SOLVING MAJOR CHALLENGES WITH A SIMPLE SOLUTION" "Простое решение основных задач"

Replace-Text "As they currently stand, both traditional cards and crypto payment solutions leave a lot to be desired. Transaction fees of over 3% for some cards can add up to substantial losses over time for businesses of all sizes. For even small businesses, it’s normal for more than `$50,000 a year to be lost in processing fees." "В настоящее время оплата традиционными картами, как и оплата криптовалютой, имеет много недостатков. Комиссия более 3% за каждую транзакцию может привести к значительным убыткам. Сейчас для малого бизнеса терять более `$50 000 в год на комиссиях является чем-то нормальным. "

Replace-Text "In addition to this, both buyers and sellers must have access to major banks and their services in order to facilitate these types of transactions; a major problem in many parts of the world. Also, traditional services inherently have long confirmation times of up to several days, which can themselves lead to other issues" "В дополнении к этому, и покупатели, и продавцы должны иметь доступ к крупным банкам и их услугам, чтобы осуществить любые переводы; в некоторых частях мира это может быть настоящей проблемой. Кроме того, использование традиционных способов оплаты занимает много времени, вплоть до нескольких дней"
